$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 1 (row 2): reward quantity changes from 15 to 5
$ws.Range("C2").Value = 5

# Day 5 (row 6): reward type changes from GOLD to CARD, quantity from 200 to 10
$ws.Range("B6").Value = "CARD"
$ws.Range("C6").Value = 10
